$wb = $excel.ActiveWorkbook

# Update both "展览" and "全部类型" sheets, which contain mirrored data.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 301
    $ws.Range("F4").Value = 172
}
